$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "DeleteProduct" worksheet after the existing "AddProduct" sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "DeleteProduct"

# Copy the header / data-row formatting from AddProduct so the new sheet
# reuses the same (bold+border / border) cell styles.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)
$ws1.Range("A2:C2").Copy()
$ws2.Range("A2:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$ws2.Range("A1").Value = "ExecutionFlag"
$ws2.Range("B1").Value = "TestCaseId"
$ws2.Range("C1").Value = "TestCaseName"
$ws2.Range("D1").Value = "ProductName"

# Test case rows for the new "Delete product from cart" scenarios
$ws2.Range("A2").Value = "Yes"
$ws2.Range("B2").Value = "TC0001"
$ws2.Range("C2").Value = "Delete Pantene Pro-V Product from Shopping Cart"
$ws2.Range("D2").Value = "Pantene Pro-V"

$ws2.Range("A3").Value = "Yes"
$ws2.Range("B3").Value = "TC0002"
$ws2.Range("C3").Value = "Delete Shaving cream Product from Shopping Cart"
$ws2.Range("D3").Value = "Shaving cream"

# Size the columns to fit their (longer) content, same as AddProduct's
# auto-fitted columns.
$ws2.Columns.Item(1).ColumnWidth = 12.592447916666666
$ws2.Columns.Item(2).ColumnWidth = 9.736979166666666
$ws2.Columns.Item(3).ColumnWidth = 46.022135416666664
$ws2.Columns.Item(4).ColumnWidth = 24.022135416666668

# Update selections to match the saved state of each sheet.
$ws1.Range("C14").Select() | Out-Null
$ws2.Range("D12").Select() | Out-Null
$ws2.Activate()

Write-Host "DeleteProduct sheet added with" ($wb.Worksheets.Count) "total sheets"
